$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a data typo: stock_quantity for product id 108 (row 9) was 5, should be 18
$ws.Range("E9").Value = 18

# Adjust row height for header + first 11 data rows (rows 1-12) from 19.5 to 18.75
for ($r = 1; $r -le 12; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
}
